# Add data for 2022-02-19 (carjacking by month YoY historical)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "through" date label
$ws.Name = "Through 2022-02-19"
$ws.Range("I1").Value = "2022 (through 02-19)"

# Update updated 2022 totals for January, February and the yearly Total row
$ws.Range("I2").Value = 159
$ws.Range("I3").Value = 92
$ws.Range("I14").Value = 251
